$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPXV_Metadata")

# Cut column AI (35) and insert it before column C (3)
$ws.Columns.Item(35).Cut() | Out-Null
$ws.Columns.Item(3).Insert() | Out-Null

# Rename the header text that used to be "ncbi_sequence_name_sra" to "ncbi-spuid-sra"
$ws.Range("C2").Value = "ncbi-spuid-sra"

# Row1 C1 should be blank (no "SRA - all" text) but keep its style
$ws.Range("C1").Value = $null
